$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Classes": update a handful of label / description cells.
# ---------------------------------------------------------------------------
$wsClasses = $wb.Worksheets.Item("Classes")

$wsClasses.Range("E2").Value = "ParteNBR"

$wsClasses.Range("F6").Value = "HospitalarSUS"
$wsClasses.Range("F7").Value = "HospitalarPRI"

$wsClasses.Range("Q6").Value = "Ambiente que pertenece a un hospital de la red del Sistema Único de Salud del Brasil"
$wsClasses.Range("Q7").Value = "Ambiente que pertenece a un hospital de la red privada del Brasil"

$wsClasses.Range("P6").Value = "Ambiente que pertence a hospital da rede do Sistema Único de Saúde do Brasil"
$wsClasses.Range("P7").Value = "Ambiente que pertence a hospital da rede privada do Brasil"

# ---------------------------------------------------------------------------
# Sheet "Proprie": columns B, E and S (rows 3-31) become formulas that copy
# the value from the row directly above, instead of repeated literal text.
# ---------------------------------------------------------------------------
$wsProprie = $wb.Worksheets.Item("Proprie")

for ($r = 3; $r -le 31; $r++) {
    $prev = $r - 1
    $wsProprie.Range("B$r").Formula = "=B$prev"
    $wsProprie.Range("E$r").Formula = "=E$prev"
    $wsProprie.Range("S$r").Formula = "=S$prev"
}

# ---------------------------------------------------------------------------
# View state: "Classes" becomes the active sheet/tab, with Q9 selected;
# "Proprie" keeps its frozen pane but the lower-left pane now has S3 selected.
# ---------------------------------------------------------------------------
$wsClasses.Select()
$wsClasses.Range("Q9").Select()

$wsProprie.Activate()
$wsProprie.Range("B22").Select()
$wsProprie.Range("S3").Select()

$wsClasses.Activate()
